# Auto-generated edit script: updates crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "27.933.54"
$ws.Range("E2").Value = "  +4.99%  "

$ws.Range("D3").Value = "1.777.90"
$ws.Range("E3").Value = "  +3.45%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "'243.37"
$ws.Range("E5").Value = "  +1.11%  "

$ws.Range("D6").Value = "'0.9998"

$ws.Range("D7").Value = "'0.4883"
$ws.Range("E7").Value = "  -0.82%  "

$ws.Range("D8").Value = "'0.2651"
$ws.Range("E8").Value = "  +1.89%  "

$ws.Range("D9").Value = "'0.06239"
$ws.Range("E9").Value = "  +0.53%  "

$ws.Range("D10").Value = "1.779.18"
$ws.Range("E10").Value = "  +3.46%  "

$ws.Range("D11").Value = "'16.29"

$ws.Range("D12").Value = "'0.07004"
$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.611"
$ws.Range("E13").Value = "  +2.92%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.6163"
$ws.Range("E14").Value = "  +1.52%  "

$ws.Range("D15").Value = "'79.39"
$ws.Range("E15").Value = "  +3.41%  "

$ws.Range("B16").Value = "Dai"
$ws.Range("C16").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  +0.22%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "27.915.21"
$ws.Range("E17").Value = "  +5.56%  "

$ws.Range("D18").Value = "'0.9992"
$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("D19").Value = "'0.000007197"
$ws.Range("E19").Value = "  +0.65%  "

$ws.Range("E20").Value = "  +3.69%  "

$ws.Range("D21").Value = "2.006.45"
$ws.Range("E21").Value = "  +2.98%  "

$ws.Range("D22").Value = "'4.561"
$ws.Range("E22").Value = "  +3.50%  "

$ws.Range("D23").Value = "'8.633"
$ws.Range("E23").Value = "  +1.43%  "

$ws.Range("D24").Value = "'5.199"
$ws.Range("E24").Value = "  +2.36%  "

$ws.Range("D25").Value = "'141.70"
$ws.Range("E25").Value = "  +2.92%  "

$ws.Range("D26").Value = "'15.58"
$ws.Range("E26").Value = "  +1.87%  "

$ws.Range("D27").Value = "'1.857"
$ws.Range("E27").Value = "  +6.53%  "

$ws.Range("D28").Value = "'108.84"
$ws.Range("E28").Value = "  +2.91%  "

$ws.Range("D30").Value = "'4.124"
$ws.Range("E30").Value = "  +5.28%  "

$ws.Range("D31").Value = "'0.08221"
$ws.Range("E31").Value = "  +3.47%  "

$ws.Range("D32").Value = "'3.777"
$ws.Range("E32").Value = "  +3.78%  "

$ws.Range("D33").Value = "'0.04747"
$ws.Range("E33").Value = "  +5.36%  "

$ws.Range("D34").Value = "'1.058"
$ws.Range("E34").Value = "  +5.98%  "

$ws.Range("D35").Value = "'2.596"
$ws.Range("E35").Value = "  -0.61%  "

$ws.Range("D36").Value = "'0.6417"
$ws.Range("E36").Value = "  +2.57%  "

$ws.Range("D37").Value = "'0.9404"
$ws.Range("E37").Value = "  +0.49%  "

$ws.Range("D38").Value = "'2.587"
$ws.Range("E38").Value = "  +7.31%  "

$ws.Range("D39").Value = "'2.037"
$ws.Range("E39").Value = "  +1.53%  "

$ws.Range("D40").Value = "'5.888"
$ws.Range("E40").Value = "  +6.81%  "

$ws.Range("E41").Value = "  +1.78%  "

$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  +0.22%  "

$ws.Range("E43").Value = "  +0.92%  "

$ws.Range("D44").Value = "'0.3939"
$ws.Range("E44").Value = "  +2.60%  "

$ws.Range("D45").Value = "'7.183"
$ws.Range("E45").Value = "  +3.80%  "

$ws.Range("D46").Value = "'0.1194"
$ws.Range("E46").Value = "  +3.49%  "

$ws.Range("E47").Value = "  +0.68%  "

$ws.Range("D48").Value = "'7.936"
$ws.Range("E48").Value = "  +2.55%  "

$ws.Range("D49").Value = "'1.281"
$ws.Range("E49").Value = "  +4.82%  "

$ws.Range("D50").Value = "'30.40"
$ws.Range("E50").Value = "  +0.90%  "

$ws.Range("D51").Value = "'52.41"
$ws.Range("E51").Value = "  +1.89%  "

